# Apply edits to the settings workbook
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# project_name: anca_panel_2 -> anca_panel_3
$ws.Range("B3").Value = "anca_panel_3"

# first_run_mode: 0 -> 1
$ws.Range("B13").Value = 1

# data_subsets: Granulos -> B
$ws.Range("B16").Value = "B"

# automatic_palette: 1 -> 0
$ws.Range("B23").Value = 0

# ccp_delta_cutoff: 0.01 -> 0.012
$ws.Range("B28").Value = 0.012

# Row 17 wraps more lines now (content length) -> row height grows
$ws.Rows.Item(17).RowHeight = 187.2

# Update selection / active cell to C14, and scroll view back to top (A4 -> A1)
[void]$ws.Range("A1").Select()
[void]$ws.Range("C14").Select()

$wb.Save()
